# Update book record for "1984" (row 2) to reflect a reservation by User_1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "Reserved"
$ws.Range("F2").Value = 44910
$ws.Range("G2").Value = 45272
$ws.Range("H2").Value = "User_1"
